$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after row 463 (before the old row 464),
# pushing all the existing data down by two rows.
$ws.Range("A464:R465").Insert()

# New row 464: "Primera" quality entry dated 2023-10-24 (serial 45223)
$ws.Cells.Item(464,1).Value2  = 11
$ws.Cells.Item(464,2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(464,3).Value2  = "Bíobío"
$ws.Cells.Item(464,4).Value2  = 45223
$ws.Cells.Item(464,5).Value2  = 8
$ws.Cells.Item(464,6).Value2  = 100112008
$ws.Cells.Item(464,7).Value2  = "Coliflor"
$ws.Cells.Item(464,8).Value2  = "Sin especificar"
$ws.Cells.Item(464,9).Value2  = "Primera"
$ws.Cells.Item(464,10).Value2 = 1000
$ws.Cells.Item(464,11).Value2 = 800
$ws.Cells.Item(464,12).Value2 = 800
$ws.Cells.Item(464,13).Value2 = 800
$ws.Cells.Item(464,14).Value2 = "`$/unidad"
$ws.Cells.Item(464,15).Value2 = "Región Metropolitana"
$ws.Cells.Item(464,16).Value2 = 800
$ws.Cells.Item(464,17).Value2 = 1
$ws.Cells.Item(464,18).Value2 = "Hortaliza"

# New row 465: "Segunda" quality entry, same date (serial 45223)
$ws.Cells.Item(465,1).Value2  = 11
$ws.Cells.Item(465,2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(465,3).Value2  = "Bíobío"
$ws.Cells.Item(465,4).Value2  = 45223
$ws.Cells.Item(465,5).Value2  = 8
$ws.Cells.Item(465,6).Value2  = 100112008
$ws.Cells.Item(465,7).Value2  = "Coliflor"
$ws.Cells.Item(465,8).Value2  = "Sin especificar"
$ws.Cells.Item(465,9).Value2  = "Segunda"
$ws.Cells.Item(465,10).Value2 = 1000
$ws.Cells.Item(465,11).Value2 = 600
$ws.Cells.Item(465,12).Value2 = 600
$ws.Cells.Item(465,13).Value2 = 600
$ws.Cells.Item(465,14).Value2 = "`$/unidad"
$ws.Cells.Item(465,15).Value2 = "Región Metropolitana"
$ws.Cells.Item(465,16).Value2 = 600
$ws.Cells.Item(465,17).Value2 = 1
$ws.Cells.Item(465,18).Value2 = "Hortaliza"
